$d = $word.ActiveDocument

# The paragraph "  SPDSit <- SPD[SPD$Site == Station, ]" is the only place in
# the document where " <- SPD[" is immediately preceded by "SPDSit" (there are
# two other "<- SPD[" occurrences involving SPDNut later on), so scope the
# Find/Replace to that specific paragraph to avoid touching the others.
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*SPDSit <- SPD[SPD`$Site*") {
        $targetPara = $p
        break
    }
}

$r = $targetPara.Range

# Merge the " <- " run with the "SPD[" run into a single run/text node, which
# also drops the now-redundant gramStart proofing marker that used to sit
# between them.
$r.Find.Execute(" <- SPD[", $true, $false, $false, $false, $false, $true, 1,
                 $false, " <- SPD[", 2)
